$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3306.5405
$ws.Range("I98").Value = 2346.8386
$ws.Range("J98").Value = 8265
$ws.Range("K98").Value = 2346.8386
$ws.Range("L98").Value = 8265
$ws.Range("M98").Value = -848.8386
$ws.Range("N98").Value = -11261

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3306.5405
$ws.Range("I122").Value = 2346.8386
$ws.Range("J122").Value = 8265
$ws.Range("K122").Value = 7040.5158
$ws.Range("L122").Value = 24795
$ws.Range("M122").Value = -4590.5158
$ws.Range("N122").Value = -29695

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2287.5
$ws.Range("I125").Value = 2341
$ws.Range("K125").Value = 21069
$ws.Range("M125").Value = -18609

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1872
$ws.Range("I127").Value = 596
$ws.Range("J127").Value = 2510
$ws.Range("K127").Value = 1788
$ws.Range("L127").Value = 7530
$ws.Range("M127").Value = 3172
$ws.Range("N127").Value = -17450

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 943.5454999999999
$ws.Range("J129").Value = 956.18604
$ws.Range("L129").Value = 2868.55812
$ws.Range("N129").Value = -12868.55812

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7927.18
$ws.Range("I32").Value = 3679.9402
$ws.Range("K32").Value = 3679.9402
$ws.Range("M32").Value = -3392.9402

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1936.4375
$ws.Range("I45").Value = 1146.1428
$ws.Range("J45").Value = 2551.111
$ws.Range("K45").Value = 1146.1428
$ws.Range("L45").Value = 2551.111
$ws.Range("M45").Value = -769.1428000000001
$ws.Range("N45").Value = -3305.111

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5962.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5962.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17886.75
$ws.Range("N132").Value = -22946.75
$ws.Range("M132").ClearContents()

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1687.2
$ws.Range("I94").Value = 1680.2222
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 1680.2222
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -1229.2222
$ws.Range("N94").Value = -2652

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3046.6667
$ws.Range("I105").Value = 3046.6667
$ws.Range("K105").Value = 3046.6667
$ws.Range("M105").Value = -1299.6667

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1456.6207
$ws.Range("I16").Value = 1147
$ws.Range("J16").Value = 1619.579
$ws.Range("K16").Value = 1147
$ws.Range("L16").Value = 1619.579
$ws.Range("M16").Value = -860
$ws.Range("N16").Value = -2193.579

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 949.125
$ws.Range("I107").Value = 599
$ws.Range("K107").Value = 599
$ws.Range("M107").Value = 1321

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1456.6207
$ws.Range("I113").Value = 1147
$ws.Range("J113").Value = 1619.579
$ws.Range("K113").Value = 1147
$ws.Range("L113").Value = 1619.579
$ws.Range("M113").Value = 1023
$ws.Range("N113").Value = -5959.579

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2184.2812
$ws.Range("I132").Value = 988.9
$ws.Range("J132").Value = 2727.6365
$ws.Range("K132").Value = 8900.1
$ws.Range("L132").Value = 24548.7285
$ws.Range("M132").Value = -6370.1
$ws.Range("N132").Value = -29608.7285

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6092.841
$ws.Range("I70").Value = 5732.2256
$ws.Range("K70").Value = 5732.2256
$ws.Range("M70").Value = -5462.2256

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6092.841
$ws.Range("I73").Value = 5732.2256
$ws.Range("K73").Value = 5732.2256
$ws.Range("M73").Value = -4796.2256

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14708630
$ws.Range("I80").Value = 35716516
$ws.Range("J80").Value = 3110.4
$ws.Range("K80").Value = 35716516
$ws.Range("L80").Value = 3110.4
$ws.Range("M80").Value = -35715518
$ws.Range("N80").Value = -5106.4

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 14708630
$ws.Range("I83").Value = 35716516
$ws.Range("J83").Value = 3110.4
$ws.Range("K83").Value = 178582580
$ws.Range("L83").Value = 15552
$ws.Range("M83").Value = -178577588
$ws.Range("N83").Value = -25536

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3023.4688
$ws.Range("I122").Value = 2617.7307
$ws.Range("K122").Value = 7853.1921
$ws.Range("M122").Value = -5403.1921

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5308.8184
$ws.Range("I132").Value = 3850
$ws.Range("J132").Value = 6142.4287
$ws.Range("K132").Value = 11550
$ws.Range("L132").Value = 18427.2861
$ws.Range("M132").Value = -9020
$ws.Range("N132").Value = -23487.2861

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6543.1816
$ws.Range("I7").Value = 5266.6665
$ws.Range("J7").Value = 7021.875
$ws.Range("K7").Value = 5266.6665
$ws.Range("L7").Value = 7021.875
$ws.Range("M7").Value = -5154.6665
$ws.Range("N7").Value = -7245.875

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1406.2
$ws.Range("I61").Value = 1292.3572
$ws.Range("K61").Value = 1292.3572
$ws.Range("M61").Value = -1090.3572

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2104.2
$ws.Range("I100").Value = 1928.8889
$ws.Range("K100").Value = 1928.8889
$ws.Range("M100").Value = -1387.8889

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1406.2
$ws.Range("I113").Value = 1292.3572
$ws.Range("K113").Value = 1292.3572
$ws.Range("M113").Value = 877.6428000000001

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3637.4285
$ws.Range("I122").Value = 3106.5625
$ws.Range("K122").Value = 9319.6875
$ws.Range("M122").Value = -6869.6875

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6543.1816
$ws.Range("I126").Value = 5266.6665
$ws.Range("J126").Value = 7021.875
$ws.Range("K126").Value = 15799.9995
$ws.Range("L126").Value = 21065.625
$ws.Range("M126").Value = -13329.9995
$ws.Range("N126").Value = -26005.625

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 29481.842
$ws.Range("J127").Value = 29481.842
$ws.Range("L127").Value = 29481.842
$ws.Range("N127").Value = -39401.842

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4515.04
$ws.Range("I132").Value = 1940.4667
$ws.Range("J132").Value = 5618.4287
$ws.Range("K132").Value = 5821.4001
$ws.Range("L132").Value = 16855.2861
$ws.Range("M132").Value = -3291.4001
$ws.Range("N132").Value = -21915.2861

# WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 31900
$ws.Range("J80").Value = 31900
$ws.Range("L80").Value = 31900
$ws.Range("N80").Value = -33896

# WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 31900
$ws.Range("J83").Value = 31900
$ws.Range("L83").Value = 95700
$ws.Range("N83").Value = -105684

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1185.8572
$ws.Range("I100").Value = 967
$ws.Range("K100").Value = 1934
$ws.Range("M100").Value = -1393

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8337893.5
$ws.Range("I132").Value = 10850.7
$ws.Range("J132").Value = 11113575
$ws.Range("K132").Value = 32552.1
$ws.Range("L132").Value = 33340725
$ws.Range("M132").Value = -30022.1
$ws.Range("N132").Value = -33345785

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3186.8235
$ws.Range("I136").Value = 1039.6316
$ws.Range("J136").Value = 5906.6
$ws.Range("K136").Value = 3118.8948
$ws.Range("L136").Value = 17719.8
$ws.Range("M136").Value = -568.8948
$ws.Range("N136").Value = -22819.8
